$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Part 1: "A Report on Buffaloes" / "By Eric Rizzi" title block becomes
# a single paragraph "A Series of Buffalo Jokes" (split over several
# runs, matching the target markup), and the blank paragraph + the old
# "By Eric Rizzi" paragraph are removed.
# ---------------------------------------------------------------------

$titleRange = $d.Content
$titleFound = $titleRange.Find.Execute("A Report on Buffaloes")
if (-not $titleFound) {
    throw "Could not find the 'A Report on Buffaloes' title paragraph"
}

# Locate the paragraph that contains the found title text so we can work
# with paragraph-relative positions (Find narrows the range to just the
# matched text, excluding the trailing paragraph mark).
$titleParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Start -le $titleRange.Start -and $cand.Range.End -ge $titleRange.End) {
        $titleParaIndex = $i
        break
    }
}
if ($titleParaIndex -eq -1) {
    throw "Could not resolve the title paragraph index"
}

$titlePara = $d.Paragraphs.Item($titleParaIndex)
$titleTextRange = $d.Range($titlePara.Range.Start, $titlePara.Range.End - 1)

# Replace the title text with the new multi-run text. InsertXML is used
# (instead of plain Range.Text / InsertAfter) so each piece of text ends
# up in its own <w:r>, matching the target markup's run layout exactly
# (the engine otherwise coalesces adjacent same-formatted runs).
$titleOoxml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:b/><w:sz w:val="40"/></w:rPr><w:t xml:space="preserve">A </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="40"/></w:rPr><w:t>Series</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="40"/></w:rPr><w:t xml:space="preserve"> o</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="40"/></w:rPr><w:t>f</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="40"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="40"/></w:rPr><w:t>Buffalo Jokes</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$titleTextRange.InsertXML($titleOoxml)

# Now remove the following blank paragraph and the "By Eric Rizzi"
# paragraph in their entirety (text + their own paragraph marks), while
# leaving the title paragraph and the paragraph after them untouched.
$blankPara = $d.Paragraphs.Item($titleParaIndex + 1)
$byEricPara = $d.Paragraphs.Item($titleParaIndex + 2)

# Delete the "By Eric Rizzi" text, keeping its paragraph mark in place
# for the moment (deleting text and marks together in one call only
# collapses one paragraph boundary, not two).
$d.Range($byEricPara.Range.Start, $byEricPara.Range.End - 1).Delete()

# Merge the (now) two adjacent empty paragraph marks away one at a time.
$d.Range($blankPara.Range.Start, $blankPara.Range.End).Delete()
$d.Range($blankPara.Range.Start, $blankPara.Range.End).Delete()

# ---------------------------------------------------------------------
# Part 2: insert a new blank paragraph (sz 56) right after "Q:  What do
# two bison in love do?", before the blank paragraph that already
# follows it.
# ---------------------------------------------------------------------

$qRange = $d.Content
$qFound = $qRange.Find.Execute("Q:  What do two bison in love do?")
if (-not $qFound) {
    throw "Could not find the 'Q:  What do two bison in love do?' paragraph"
}

$qParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Start -le $qRange.Start -and $cand.Range.End -ge $qRange.End) {
        $qParaIndex = $i
        break
    }
}
if ($qParaIndex -eq -1) {
    throw "Could not resolve the 'two bison in love' paragraph index"
}

$qPara = $d.Paragraphs.Item($qParaIndex)
$insertPoint = $d.Range($qPara.Range.End, $qPara.Range.End)

$newParaOoxml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:sz w:val="56"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$insertPoint.InsertXML($newParaOoxml)

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
